# Insert a new weekly observation row into the Espinaca (Femacal de La Calera)
# dataset. The new record is inserted at row 194, pushing the existing rows
# 194:245 down to 195:246 (row 245's former data ends up in row 246).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 194 - this shifts rows 194:245 down to 195:246
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new record's data
$ws.Cells.Item(194, 1).Value = 3
$ws.Cells.Item(194, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(194, 3).Value = "Coquimbo"
$ws.Cells.Item(194, 4).Value = 44551
$ws.Cells.Item(194, 5).Value = 5
$ws.Cells.Item(194, 6).Value = 100112012
$ws.Cells.Item(194, 7).Value = "Espinaca"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 165
$ws.Cells.Item(194, 11).Value = 4500
$ws.Cells.Item(194, 12).Value = 5000
$ws.Cells.Item(194, 13).Value = 4742
$ws.Cells.Item(194, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(194, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(194, 16).Value = 1581
$ws.Cells.Item(194, 17).Value = 3
$ws.Cells.Item(194, 18).Value = "Hortaliza"

# Give the date cell the same date/time number format used by the rest of
# column D (style carried over automatically from the row-insert, but set
# explicitly here for safety).
$ws.Cells.Item(194, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
